$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column J (10), shifting the old J column (完成日期) to M
$ws.Range("J1:L1").EntireColumn.Insert()

# Set the new header cell values
$ws.Range("J2").Value = "故障說明"
$ws.Range("K2").Value = "維修經過"
$ws.Range("L2").Value = "結果"

# Give the three newly inserted columns a width; column M keeps the width
# inherited from the old column J (14.375) automatically via the insert above.
$ws.Range("J1:L1").EntireColumn.ColumnWidth = 10.71

Write-Host "Done"
